$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the registration success messages (en/ru) to mention e-mail activation
$ws.Range("C2").Value = "User registered successfully. Check your e-mail to complete registration."
$ws.Range("D2").Value = "Пользователь успешно зарегистрирован. Письмо с подтверждением регистрации от правлено вам на почту."

# The longer text now wraps onto a 3rd line, so the row needs to grow accordingly
$ws.Rows.Item(2).RowHeight = 45

# Reflect where the user left the selection when saving
$ws.Range("D5").Select()
